# Updates the cryptos price table to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: safe to assign directly.
$textUpdates = @{
    'D2' = '30.163.18'
    'E2' = '  +4.81%  '
    'D3' = '1.915.02'
    'E3' = '  +6.05%  '
    'E4' = '  -0.11%  '
    'E5' = '  +1.45%  '
    'E6' = '  -0.11%  '
    'E7' = '  +5.12%  '
    'B8' = 'OKB'
    'C8' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'E8' = '  +6.87%  '
    'B9' = 'Cardano'
    'C9' = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
    'E9' = '  +7.39%  '
    'B10' = 'Dogecoin'
    'C10' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'E10' = '  +6.53%  '
    'B11' = 'WrappedEther'
    'C11' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D11' = '1.914.28'
    'E11' = '  +6.06%  '
    'B12' = 'Solana'
    'C12' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'E12' = '  +3.70%  '
    'B13' = 'TRON'
    'C13' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'E13' = '  +3.32%  '
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E14' = '  +7.44%  '
    'B15' = 'Litecoin'
    'C15' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'E15' = '  +7.67%  '
    'B16' = 'Polkadot'
    'C16' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'E16' = '  +4.67%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '30.159.71'
    'E17' = '  +4.90%  '
    'B18' = 'ShibaInu'
    'C18' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'E18' = '  +6.23%  '
    'B19' = 'Dai'
    'C19' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'E19' = '  +0.09%  '
    'B20' = 'Avalanche'
    'C20' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'E20' = '  +7.29%  '
    'B21' = 'WrappedliquidstakedEther2.0'
    'C21' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D21' = '2.157.51'
    'E21' = '  +5.96%  '
    'B22' = 'BinanceUSD'
    'C22' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'E22' = '  -0.17%  '
    'B23' = 'Uniswap'
    'C23' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'E23' = '  +5.91%  '
    'B24' = 'Chainlink'
    'C24' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E24' = '  +8.23%  '
    'B25' = 'Cosmos'
    'C25' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E25' = '  +4.12%  '
    'B26' = 'BitcoinCash'
    'C26' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'E26' = '  +26.30%  '
    'B27' = 'Monero'
    'C27' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E27' = '  +2.71%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'E28' = '  +8.07%  '
    'B29' = 'LidoDAOToken'
    'C29' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'E29' = '  +7.32%  '
    'B30' = 'Toncoin'
    'C30' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'E30' = '  -1.15%  '
    'B31' = 'InternetComputer(DFINITY)'
    'C31' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E31' = '  +2.65%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'E32' = '  +5.92%  '
    'B33' = 'Filecoin'
    'C33' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E33' = '  +6.07%  '
    'B34' = 'Hedera'
    'C34' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E34' = '  +3.69%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'E35' = '  +6.50%  '
    'B36' = 'ImmutableX'
    'C36' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'E36' = '  +8.59%  '
    'B37' = 'HuobiToken'
    'C37' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'E37' = '  +0.61%  '
    'E38' = '  +7.99%  '
    'B39' = 'RenderToken'
    'C39' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'E39' = '  +7.61%  '
    'B40' = 'TrustWalletToken'
    'C40' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E40' = '  +2.24%  '
    'B41' = 'VeChain'
    'C41' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E41' = '  +5.67%  '
    'B42' = 'FraxShare'
    'C42' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E42' = '  +3.68%  '
    'E43' = '  +6.04%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E44' = '  +5.49%  '
    'B45' = 'PaxDollar'
    'C45' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'E45' = '  +0.02%  '
    'B46' = 'Aptos'
    'C46' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'E46' = '  +6.84%  '
    'B47' = 'Algorand'
    'C47' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E47' = '  +4.40%  '
    'B48' = 'Cronos'
    'C48' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'E48' = '  +4.31%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E49' = '  +5.22%  '
    'B50' = 'Elrond'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'E50' = '  +6.25%  '
    'B51' = 'Decentraland'
    'C51' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'E51' = '  +6.77%  '
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

# Numeric-looking values (e.g. "88.00", "0.9990") must be forced to stay as text
# so Excel does not silently convert them to doubles and lose formatting/leading
# zeros or the multi-dot "thousands" notation used by the source site.
$numericTextUpdates = @{
    'D4' = '0.9983'
    'D5' = '252.67'
    'D6' = '0.9983'
    'D7' = '0.5197'
    'D8' = '46.31'
    'D9' = '0.2995'
    'D10' = '0.06774'
    'D12' = '17.59'
    'D13' = '0.07316'
    'D14' = '0.6915'
    'D15' = '88.00'
    'D16' = '4.916'
    'D18' = '0.000007791'
    'D19' = '0.9992'
    'D20' = '13.16'
    'D22' = '0.9977'
    'D23' = '4.873'
    'D24' = '5.741'
    'D25' = '9.240'
    'D26' = '140.30'
    'D27' = '146.49'
    'D28' = '17.26'
    'D29' = '2.017'
    'D30' = '1.376'
    'D31' = '4.285'
    'D32' = '0.08851'
    'D33' = '4.066'
    'D34' = '0.05153'
    'D35' = '1.160'
    'D36' = '0.7262'
    'D37' = '2.693'
    'D38' = '2.836'
    'D39' = '2.303'
    'D40' = '0.9708'
    'D41' = '0.01689'
    'D42' = '6.189'
    'D43' = '0.4352'
    'D44' = '106.19'
    'D45' = '0.9990'
    'D46' = '7.710'
    'D47' = '0.1276'
    'D48' = '0.05722'
    'D49' = '8.564'
    'D50' = '33.15'
    'D51' = '0.3861'
}

foreach ($cellRef in $numericTextUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$cellRef]
    $cell.ClearFormats()
}
